$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 14 data: SS-11
$ws.Range("A14").Value = "SS-11"
$ws.Range("B14").Value = 9.32
$ws.Range("C14").Value = 9.32
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = 1.8
$ws.Range("F14").Value = 84
$ws.Range("G14").Value = 15
$ws.Range("N14").Value = 100

# New row 15 data: SS-12
$ws.Range("A15").Value = "SS-12"
$ws.Range("B15").Value = 9.32
$ws.Range("C15").Value = 10.51
$ws.Range("D15").Value = 10
$ws.Range("E15").Value = 1.8
$ws.Range("F15").Value = 84
$ws.Range("G15").Value = 15
$ws.Range("N15").Value = 100

# Update the selected cell to match the new active cell
$ws.Range("G14").Select()
